$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "last updated" timestamp string (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 26 de Agosto de 2020 a las 09:06"

# --- Row 4: Estados Unidos ---
$ws.Range("B4").Value = 5956036
$ws.Range("C4").Value = 308
$ws.Range("D4").Value = 3254739
$ws.Range("E4").Value = 2518885
$ws.Range("G4").Value = 8
$ws.Range("H4").Value = 182412

# --- Row 6: India ---
$ws.Range("B6").Value = 3235725
$ws.Range("C6").Value = 3971
$ws.Range("D6").Value = 2468688
$ws.Range("E6").Value = 707405
$ws.Range("G6").Value = 20
$ws.Range("H6").Value = 59632

# --- Rows 30/31: Ucrania overtakes Ecuador in ranking (rows swap) ---
# Row 30 becomes Ucrania with refreshed totals
$ws.Range("A30").Value = "Ucrania"
$ws.Range("B30").Value = 110085
$ws.Range("C30").Value = 1670
$ws.Range("D30").Value = 53454
$ws.Range("E30").Value = 54313
$ws.Range("F30").Value = 0
$ws.Range("G30").Value = 0
$ws.Range("H30").Value = 2318

# Row 31 becomes Ecuador with its prior (unchanged) totals
$ws.Range("A31").Value = "Ecuador"
$ws.Range("B31").Value = 109030
$ws.Range("C31").Value = 0
$ws.Range("D31").Value = 95025
$ws.Range("E31").Value = 7637
$ws.Range("F31").Value = 0
$ws.Range("G31").Value = 0
$ws.Range("H31").Value = 6368

# --- Row 58: Armenia ---
$ws.Range("B58").Value = 43067
$ws.Range("C58").Value = 131
$ws.Range("D58").Value = 36726
$ws.Range("E58").Value = 5480
$ws.Range("G58").Value = 3
$ws.Range("H58").Value = 861

# --- Row 61: Suiza ---
$ws.Range("D61").Value = 34800
$ws.Range("E61").Value = 3460

# --- Row 72: Australia ---
$ws.Range("D72").Value = 20100
$ws.Range("E72").Value = 4555

# --- Row 109: Hungria ---
$ws.Range("B109").Value = 5288
$ws.Range("C109").Value = 73
$ws.Range("D109").Value = 3734
$ws.Range("E109").Value = 940

# --- Row 111: Hong Kong ---
$ws.Range("E111").Value = 525
$ws.Range("G111").Value = 1
$ws.Range("H111").Value = 78

# --- Row 151: Georgia ---
$ws.Range("B151").Value = 1436
$ws.Range("C151").Value = 7
$ws.Range("E151").Value = 267

# --- Row 171: Taiwan ---
$ws.Range("D171").Value = 462
$ws.Range("E171").Value = 18
